$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Remove the bold/centered/bordered header style by resetting the header
# row cells back to the default (General) style.
$ws.Range("A1:Q1").Style = "Normal"

# Update row 2 values
$ws.Range("B2").Value = 0.72
$ws.Range("C2").Value = 0.5126185691505004
$ws.Range("G2").Value = 0.05711068686699131
$ws.Range("H2").Value = 4.248052843601034

# Update row 3 values
$ws.Range("A3").Value = 492
$ws.Range("B3").Value = 22.24
$ws.Range("C3").Value = 2.878962123794901
$ws.Range("E3").Value = 22.63584059186897
$ws.Range("G3").Value = 1.772538660352875
$ws.Range("H3").Value = 66.68663515389257
$ws.Range("J3").Value = 1.389304899039429
$ws.Range("L3").Value = 321.86

# Update row 4 values
$ws.Range("A4").Value = 484
$ws.Range("B4").Value = 50.19
$ws.Range("C4").Value = 4.360260466518791
$ws.Range("D4").Value = 4.36
$ws.Range("E4").Value = 36.0926112263078
$ws.Range("F4").Value = 36.02
$ws.Range("G4").Value = 2.847068023443553
$ws.Range("H4").Value = 127.7211490456516
$ws.Range("J4").Value = 2.660857271784409
$ws.Range("L4").Value = 918.85
